$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("C1").Value = "qty"
$ws.Range("D1").Value = "startDate"
$ws.Range("E1").Value = "finishDate"

# Update data row values (dates as text)
$ws.Range("D2").Value = "10/01/2019"
$ws.Range("E2").Value = "10/26/2019"

# Remove now-unused trailing columns F:I for both rows
$ws.Range("F1:I2").ClearContents()

# Update the selected/active cell as recorded in the saved view state
$ws.Range("I8").Select()
